$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

$ws1.Range("C2").Value = 0.646323325211804
$ws1.Range("D2").Value = 0.353901374928943
$ws1.Range("E2").Value = 1.00022470014075
$ws1.Range("J2").Value = 0.353821871104706
$ws1.Range("K2").Value = 0.230665882890985
$ws1.Range("L2").Value = 0.0141349591889687
$ws1.Range("M2").Value = 0.0786918603705379
$ws1.Range("N2").Value = 0.244800842079954
$ws1.Range("F3").Value = 0.581751918085453
$ws1.Range("G3").Value = 0.230717713547336
$ws1.Range("H4").Value = 0.567613782769165
$ws1.Range("I4").Value = 0.210720923040043
$ws1.Range("O4").Value = 0.432513731475244
$ws1.Range("C5").Value = 0.87736090521989
$ws1.Range("D5").Value = 0.122717006372813
$ws1.Range("E5").Value = 1.0000779115927
$ws1.Range("J5").Value = 0.122707446040255
$ws1.Range("K5").Value = 0.0557511374272632
$ws1.Range("L5").Value = 0.0167761595032997
$ws1.Range("M5").Value = 0.0722841742237455
$ws1.Range("N5").Value = 0.0725272969305629
$ws1.Range("F6").Value = 0.82184856578161
$ws1.Range("G6").Value = 0.0557554810871752
$ws1.Range("H7").Value = 0.805071099221004
$ws1.Range("I7").Value = 0.040610865929157
$ws1.Range("O7").Value = 0.194991620264001
$ws1.Range("C8").Value = 0.835673869589237
$ws1.Range("D8").Value = 0.164430530342826
$ws1.Range("E8").Value = 1.00010439993206
$ws1.Range("J8").Value = 0.164413365598627
$ws1.Range("K8").Value = 0.106566271792464
$ws1.Range("L8").Value = 0.0481217115529761
$ws1.Range("M8").Value = 0.0557535126237842
$ws1.Range("N8").Value = 0.15468798334544
$ws1.Range("F9").Value = 0.828041271758916
$ws1.Range("G9").Value = 0.106577397304
$ws1.Range("H10").Value = 0.779914536302523
$ws1.Range("I10").Value = 0.0935344699065959
$ws1.Range("O10").Value = 0.220166878222411
$ws1.Range("C11").Value = 0.768779677316853
$ws1.Range("D11").Value = 0.231367222074937
$ws1.Range("E11").Value = 1.00014689939179
$ws1.Range("J11").Value = 0.231333239362774
$ws1.Range("K11").Value = 0.186155599695415
$ws1.Range("L11").Value = 0.0173592097369919
$ws1.Range("M11").Value = 0.00721785739467148
$ws1.Range("N11").Value = 0.203514809432407
$ws1.Range("F12").Value = 0.778922519417664
$ws1.Range("G12").Value = 0.186182945839789
$ws1.Range("H13").Value = 0.76156075962332
$ws1.Range("I13").Value = 0.171101084947779
$ws1.Range("O13").Value = 0.238551096757445
$ws1.Range("C14").Value = 0.86056319164205
$ws1.Range("D14").Value = 0.139525364609732
$ws1.Range("E14").Value = 1.00008855625178
$ws1.Range("J14").Value = 0.139513009860504
$ws1.Range("K14").Value = 0.125378383607537
$ws1.Range("L14").Value = 0.00832235281708016
$ws1.Range("M14").Value = 0.0356365819701996
$ws1.Range("N14").Value = 0.133700736424617
$ws1.Range("F15").Value = 0.833246543643177
$ws1.Range("G15").Value = 0.125389486647244
$ws1.Range("H16").Value = 0.824923453829725
$ws1.Range("I16").Value = 0.134747448835763
$ws1.Range("O16").Value = 0.175149591830703

# --- Sheet 2: "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

$ws2.Range("B2").Value = 0.232938856719777
$ws2.Range("C2").Value = 0.256662827440131
$ws2.Range("D2").Value = 0.244800842079954
$ws2.Range("B3").Value = 0.405015812137435
$ws2.Range("C3").Value = 0.460011650813053
$ws2.Range("D3").Value = 0.432513731475244
$ws2.Range("B4").Value = 0.305652492348059
$ws2.Range("C4").Value = 0.401991249861353
$ws2.Range("D4").Value = 0.353821871104706
$ws2.Range("B5").Value = 0.0690421811631757
$ws2.Range("C5").Value = 0.0760124126979502
$ws2.Range("D5").Value = 0.0725272969305629
$ws2.Range("B6").Value = 0.156308979164295
$ws2.Range("C6").Value = 0.233674261363707
$ws2.Range("D6").Value = 0.194991620264001
$ws2.Range("B7").Value = 0.06715017383804
$ws2.Range("C7").Value = 0.178264718242471
$ws2.Range("D7").Value = 0.122707446040255
$ws2.Range("B8").Value = 0.126334894574742
$ws2.Range("C8").Value = 0.183041072116139
$ws2.Range("D8").Value = 0.15468798334544
$ws2.Range("B9").Value = 0.0772296102465609
$ws2.Range("C9").Value = 0.363104146198261
$ws2.Range("D9").Value = 0.220166878222411
$ws2.Range("B10").Value = 0.0513894750191
$ws2.Range("C10").Value = 0.277437256178154
$ws2.Range("D10").Value = 0.164413365598627
$ws2.Range("B11").Value = 0.136666777129534
$ws2.Range("C11").Value = 0.27036284173528
$ws2.Range("D11").Value = 0.203514809432407
$ws2.Range("B12").Value = -0.0115602392328289
$ws2.Range("C12").Value = 0.488662432747719
$ws2.Range("D12").Value = 0.238551096757445
$ws2.Range("B13").Value = 0.0781908552653348
$ws2.Range("C13").Value = 0.384475623460212
$ws2.Range("D13").Value = 0.231333239362774
$ws2.Range("B14").Value = 0.123211014591888
$ws2.Range("C14").Value = 0.144190458257347
$ws2.Range("D14").Value = 0.133700736424617
$ws2.Range("B15").Value = 0.110434530667312
$ws2.Range("C15").Value = 0.239864652994094
$ws2.Range("D15").Value = 0.175149591830703
$ws2.Range("B16").Value = 0.0744454696015892
$ws2.Range("C16").Value = 0.204580550119418
$ws2.Range("D16").Value = 0.139513009860504
